$d = $word.ActiveDocument

# The "Preparation" section currently contains two numbered practice
# questions (custom numbering, numId 1001 -> wdListSimpleNumbering) and a
# trailing "-->" marker paragraph right before the "Solutions" bookmark.
# All three paragraphs need to be removed, leaving "Preparation" followed
# directly by the Solutions section.

$paras = @($d.Paragraphs)

$numberedIdx = @()
for ($i = 0; $i -lt $paras.Count; $i++) {
    if ($paras[$i].Range.ListFormat.ListType -eq 4) {
        $numberedIdx += $i
    }
}

if ($numberedIdx.Count -eq 0) {
    # Fallback: locate the questions by their known text if list-type
    # detection doesn't find anything (defensive; shouldn't normally hit).
    for ($i = 0; $i -lt $paras.Count; $i++) {
        $t = $paras[$i].Range.Text.Trim()
        if ($t.StartsWith("An eight sided die") -or $t.StartsWith("Calculate the probability of not getting")) {
            $numberedIdx += $i
        }
    }
}

if ($numberedIdx.Count -gt 0) {
    $firstIdx = $numberedIdx[0]
    $lastIdx = $numberedIdx[$numberedIdx.Count - 1]

    # Also sweep up the "-->" marker paragraph that immediately follows
    # the numbered block, if present.
    $endIdx = $lastIdx
    if (($lastIdx + 1) -lt $paras.Count) {
        $nextText = $paras[$lastIdx + 1].Range.Text.Trim()
        if ($nextText -eq "–>") {
            $endIdx = $lastIdx + 1
        }
    }

    $startPos = $paras[$firstIdx].Range.Start
    $endPos = $paras[$endIdx].Range.End
    $rng = $d.Range($startPos, $endPos)
    $rng.Delete()
}
